# Update the localization status from "Ready for handoff" to "In Translation"
# across all three worksheets (Overview, zh-cn, de-de), then shrink the
# status columns to match the narrower text (mirrors Excel's own column
# autosize once the shorter "In Translation" label is entered).

$wb = $excel.ActiveWorkbook

# --- Overview sheet: status columns are E (zh-cn) and F (de-de), rows 2-4 ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F4").Value = "In Translation"
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# --- zh-cn sheet: Status is column C, rows 2-4 ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2:C4").Value = "In Translation"
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

# --- de-de sheet: Status is column C, rows 2-4 ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2:C4").Value = "In Translation"
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
